$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rescale B/D/H columns (rows 1-4) by 1e-15 and format as scientific notation
$cols = @("B","D","H")
for ($row = 1; $row -le 4; $row++) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $rng = $ws.Range($addr)
        $v = $rng.Value2
        $rng.Value2 = ($v / 100) / 10000000000000
        $rng.NumberFormat = "0.00E+00"
    }
}

# Update sheet view: clear the scrolled topLeftCell and move selection to D5
$ws.Range("D5").Select()
